$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new Wins/Losses/Ties columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, bordered, centered) from an existing header cell (AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record values for every data row (2-50): 81 wins, 81 losses, 0 ties
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 81
    $ws.Cells.Item($r, 31).Value = 81
    $ws.Cells.Item($r, 32).Value = 0
}
